$wb = $excel.ActiveWorkbook
${ws1} = $wb.Worksheets.Item(1)
${ws2} = $wb.Worksheets.Item(2)
${ws3} = $wb.Worksheets.Item(3)
${ws4} = $wb.Worksheets.Item(4)
${ws1}.Range("F3").Value2 = 1736
${ws1}.Range("F4").Value2 = 108
${ws1}.Range("F5").Value2 = 72
${ws1}.Range("F6").Value2 = 711
${ws1}.Range("F7").Value2 = 177
${ws1}.Range("F8").Value2 = 213
${ws1}.Range("F9").Value2 = 19
${ws1}.Range("F10").Value2 = 1310
${ws1}.Range("F11").Value2 = 36
${ws1}.Range("F12").Value2 = 569
${ws1}.Range("F13").Value2 = 486
${ws1}.Range("F14").Value2 = 134
${ws1}.Range("F17").Value2 = 784
${ws1}.Range("F18").Value2 = 2584
${ws1}.Range("F21").Value2 = 8
${ws1}.Range("F23").Value2 = 191
${ws1}.Range("F25").Value2 = 137
${ws1}.Range("F27").Value2 = 948
${ws1}.Range("F29").Value2 = 172
${ws1}.Range("F32").Value2 = 36
${ws1}.Range("F33").Value2 = 258
${ws2}.Range("F4").Value2 = 627
${ws2}.Range("F5").Value2 = 627
${ws2}.Range("F10").Value2 = 294
${ws2}.Range("F14").Value2 = 535
${ws2}.Range("F29").Value2 = 191
${ws3}.Range("F5").Value2 = 2360
${ws3}.Range("F6").Value2 = 953
${ws3}.Range("F9").Value2 = 1196
${ws4}.Range("F4").Value2 = 2360
${ws4}.Range("F6").Value2 = 1736
${ws4}.Range("F8").Value2 = 953
${ws4}.Range("F9").Value2 = 1196
${ws4}.Range("F12").Value2 = 109
${ws4}.Range("F13").Value2 = 72
${ws4}.Range("F14").Value2 = 711
${ws4}.Range("F15").Value2 = 177
${ws4}.Range("F17").Value2 = 214
${ws4}.Range("F18").Value2 = 19
${ws4}.Range("F19").Value2 = 36
${ws4}.Range("F20").Value2 = 569
${ws4}.Range("F21").Value2 = 627
${ws4}.Range("F22").Value2 = 486
${ws4}.Range("F23").Value2 = 134
${ws4}.Range("F26").Value2 = 784
${ws4}.Range("F27").Value2 = 2584
${ws4}.Range("F31").Value2 = 191
${ws4}.Range("F32").Value2 = 137
${ws4}.Range("F34").Value2 = 948
${ws4}.Range("F35").Value2 = 535
${ws4}.Range("F38").Value2 = 172
${ws4}.Range("F47").Value2 = 191
${ws4}.Range("F49").Value2 = 36
${ws4}.Range("F50").Value2 = 258
